$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused duplicate columns U:AD (data that used to repeat
# columns K:T under a second copy of the K..T headers).
$ws.Range("U1:AD23").Delete()

# Rebuild the A1:T23 block with the updated Holden-scheme simulation data.
$arr = New-Object 'object[,]' 23,20
$arr[0,1] = 0
$arr[0,2] = 1
$arr[0,3] = 2
$arr[0,4] = 3
$arr[0,5] = 4
$arr[0,6] = 5
$arr[0,7] = 6
$arr[0,8] = 7
$arr[0,9] = 8
$arr[0,10] = 9
$arr[0,11] = 10
$arr[0,12] = 11
$arr[0,13] = 12
$arr[0,14] = 13
$arr[0,15] = 14
$arr[0,16] = 15
$arr[0,17] = 16
$arr[0,18] = 17
$arr[0,19] = 18
$arr[1,0] = 0
$arr[1,1] = 'HKL'
$arr[1,2] = '[3, 2, 1]'
$arr[1,3] = '[3, 1, 0]'
$arr[1,4] = '[2, 2, 2]'
$arr[1,5] = '[1, 1, 0]'
$arr[1,6] = '[2, 0, 0]'
$arr[1,7] = '[2, 2, 0]'
$arr[1,8] = '[4, 0, 0]'
$arr[1,9] = '[2, 1, 1]'
$arr[1,10] = '1Pair-A'
$arr[1,11] = '1Pair-B'
$arr[1,12] = '2Pairs-A'
$arr[1,13] = '2Pairs-B'
$arr[1,14] = '3Pairs-A'
$arr[1,15] = '3Pairs-B'
$arr[1,16] = '3Pairs-C'
$arr[1,17] = '4Pairs'
$arr[1,18] = '5A4F'
$arr[1,19] = 'MaxUnique'
$arr[2,0] = 1
$arr[2,1] = 'BT8Hex_2.5'
$arr[2,2] = 1.000016338396526
$arr[2,3] = 0.9999640576297208
$arr[2,4] = 1.000043567909944
$arr[2,5] = 1.000016338396526
$arr[2,6] = 0.9999346473816181
$arr[2,7] = 1.000016338396526
$arr[2,8] = 0.9999346473816181
$arr[2,9] = 1.000016338396526
$arr[2,10] = 1.000016338396526
$arr[2,11] = 1.000016338396526
$arr[2,12] = 0.9999754928890721
$arr[2,13] = 0.9999754928890721
$arr[2,14] = 0.999971681135955
$arr[2,15] = 0.9999891080582234
$arr[2,16] = 0.9999891080582234
$arr[2,17] = 0.9999959156427991
$arr[2,18] = 0.9999959156427991
$arr[2,19] = 0.999998548018477
$arr[3,0] = 2
$arr[3,1] = 'BT8Hex_5'
$arr[3,2] = 1.000031539742068
$arr[3,3] = 0.9999306125977502
$arr[3,4] = 1.000084107185589
$arr[3,5] = 1.000031539742068
$arr[3,6] = 0.9998738379158184
$arr[3,7] = 1.000031539742068
$arr[3,8] = 0.9998738379158184
$arr[3,9] = 1.000031539742068
$arr[3,10] = 1.000031539742068
$arr[3,11] = 1.000031539742068
$arr[3,12] = 0.9999526888289432
$arr[3,13] = 0.9999526888289432
$arr[3,14] = 0.9999453300852122
$arr[3,15] = 0.9999789724666514
$arr[3,16] = 0.9999789724666514
$arr[3,17] = 0.9999921142855055
$arr[3,18] = 0.9999921142855055
$arr[3,19] = 0.999997196154227
$arr[4,0] = 3
$arr[4,1] = 'BT8Hex_10'
$arr[4,2] = 1.000060441175615
$arr[4,3] = 0.9998670266168441
$arr[4,4] = 1.000161184662702
$arr[4,5] = 1.000060441175615
$arr[4,6] = 0.9997582249966338
$arr[4,7] = 1.000060441175615
$arr[4,8] = 0.9997582249966338
$arr[4,9] = 1.000060441175615
$arr[4,10] = 1.000060441175615
$arr[4,11] = 1.000060441175615
$arr[4,12] = 0.9999093330861244
$arr[4,13] = 0.9999093330861244
$arr[4,14] = 0.9998952309296977
$arr[4,15] = 0.999959702449288
$arr[4,16] = 0.999959702449288
$arr[4,17] = 0.9999848871308697
$arr[4,18] = 0.9999848871308697
$arr[4,19] = 0.9999946266338376
$arr[5,0] = 4
$arr[5,1] = 'BT8Hex_15'
$arr[5,2] = 1.00008864396666
$arr[5,3] = 0.9998049717797023
$arr[5,4] = 1.000236398302213
$arr[5,5] = 1.00008864396666
$arr[5,6] = 0.9996454002629926
$arr[5,7] = 1.00008864396666
$arr[5,8] = 0.9996454002629926
$arr[5,9] = 1.00008864396666
$arr[5,10] = 1.00008864396666
$arr[5,11] = 1.00008864396666
$arr[5,12] = 0.9998670221148264
$arr[5,13] = 0.9998670221148264
$arr[5,14] = 0.999846338669785
$arr[5,15] = 0.9999408960654376
$arr[5,16] = 0.9999408960654376
$arr[5,17] = 0.9999778330407432
$arr[5,18] = 0.9999778330407432
$arr[5,19] = 0.9999921170408146
$arr[6,0] = 5
$arr[6,1] = 'Spiral2.5'
$arr[6,2] = 1.000001281691743
$arr[6,3] = 0.9999971839206597
$arr[6,4] = 1.000003414213624
$arr[6,5] = 1.000001281691743
$arr[6,6] = 0.9999948771579318
$arr[6,7] = 1.000001281691743
$arr[6,8] = 0.9999948771579318
$arr[6,9] = 1.000001281691743
$arr[6,10] = 1.000001281691743
$arr[6,11] = 1.000001281691743
$arr[6,12] = 0.9999980794248375
$arr[6,13] = 0.9999980794248375
$arr[6,14] = 0.9999977809234449
$arr[6,15] = 0.9999991468471393
$arr[6,16] = 0.9999991468471393
$arr[6,17] = 0.9999996805582902
$arr[6,18] = 0.9999996805582902
$arr[6,19] = 0.9999998867279074
$arr[7,0] = 6
$arr[7,1] = 'Spiral5'
$arr[7,2] = 1.000002899690642
$arr[7,3] = 0.9999936241919529
$arr[7,4] = 1.00000772897914
$arr[7,5] = 1.000002899690642
$arr[7,6] = 0.9999884049616453
$arr[7,7] = 1.000002899690642
$arr[7,8] = 0.9999884049616453
$arr[7,9] = 1.000002899690642
$arr[7,10] = 1.000002899690642
$arr[7,11] = 1.000002899690642
$arr[7,12] = 0.9999956523261437
$arr[7,13] = 0.9999956523261437
$arr[7,14] = 0.9999949762814134
$arr[7,15] = 0.9999980681143098
$arr[7,16] = 0.9999980681143098
$arr[7,17] = 0.9999992760083929
$arr[7,18] = 0.9999992760083929
$arr[7,19] = 0.9999997428674442
$arr[8,0] = 7
$arr[8,1] = 'Spiral7.5'
$arr[8,2] = 1.000004960130395
$arr[8,3] = 0.9999890924653462
$arr[8,4] = 1.000013221860407
$arr[8,5] = 1.000004960130395
$arr[8,6] = 0.9999801647846362
$arr[8,7] = 1.000004960130395
$arr[8,8] = 0.9999801647846362
$arr[8,9] = 1.000004960130395
$arr[8,10] = 1.000004960130395
$arr[8,11] = 1.000004960130395
$arr[8,12] = 0.9999925624575157
$arr[8,13] = 0.9999925624575157
$arr[8,14] = 0.9999914057934592
$arr[8,15] = 0.9999966950151421
$arr[8,16] = 0.9999966950151421
$arr[8,17] = 0.9999987612939554
$arr[8,18] = 0.9999987612939554
$arr[8,19] = 0.999999559916929
$arr[9,0] = 8
$arr[9,1] = 'Spiral10'
$arr[9,2] = 1.000009950724764
$arr[9,3] = 0.9999781122502857
$arr[9,4] = 1.000026532059463
$arr[9,5] = 1.000009950724764
$arr[9,6] = 0.99996020007389
$arr[9,7] = 1.000009950724764
$arr[9,8] = 0.99996020007389
$arr[9,9] = 1.000009950724764
$arr[9,10] = 1.000009950724764
$arr[9,11] = 1.000009950724764
$arr[9,12] = 0.9999850753993267
$arr[9,13] = 0.9999850753993267
$arr[9,14] = 0.9999827543496463
$arr[9,15] = 0.9999933671744724
$arr[9,16] = 0.9999933671744724
$arr[9,17] = 0.9999975130620451
$arr[9,18] = 0.9999975130620451
$arr[9,19] = 0.9999991160929881
$arr[10,0] = 9
$arr[10,1] = 'Spiral15'
$arr[10,2] = 1.000017947384466
$arr[10,3] = 0.9999605203613787
$arr[10,4] = 1.000047852440248
$arr[10,5] = 1.000017947384466
$arr[10,6] = 0.9999282168972887
$arr[10,7] = 1.000017947384466
$arr[10,8] = 0.9999282168972887
$arr[10,9] = 1.000017947384466
$arr[10,10] = 1.000017947384466
$arr[10,11] = 1.000017947384466
$arr[10,12] = 0.9999730821408772
$arr[10,13] = 0.9999730821408772
$arr[10,14] = 0.9999688948810443
$arr[10,15] = 0.9999880372220734
$arr[10,16] = 0.9999880372220734
$arr[10,17] = 0.9999955147626715
$arr[10,18] = 0.9999955147626715
$arr[10,19] = 0.9999984053087188
$arr[11,0] = 10
$arr[11,1] = 'OffsetF45'
$arr[11,2] = 1.000028097300705
$arr[11,3] = 0.999938188232406
$arr[11,4] = 1.000074918695383
$arr[11,5] = 1.000028097300705
$arr[11,6] = 0.9998876122407555
$arr[11,7] = 1.000028097300705
$arr[11,8] = 0.9998876122407555
$arr[11,9] = 1.000028097300705
$arr[11,10] = 1.000028097300705
$arr[11,11] = 1.000028097300705
$arr[11,12] = 0.9999578547707304
$arr[11,13] = 0.9999578547707304
$arr[11,14] = 0.9999512992579556
$arr[11,15] = 0.9999812689473887
$arr[11,16] = 0.9999812689473887
$arr[11,17] = 0.9999929760357178
$arr[11,18] = 0.9999929760357178
$arr[11,19] = 0.9999975018451099
$arr[12,0] = 11
$arr[12,1] = 'OffsetA45'
$arr[12,2] = 1.000007626290974
$arr[12,3] = 0.9999832304387672
$arr[12,4] = 1.000020322947069
$arr[12,5] = 1.000007626290974
$arr[12,6] = 0.999969510855254
$arr[12,7] = 1.000007626290974
$arr[12,8] = 0.999969510855254
$arr[12,9] = 1.000007626290974
$arr[12,10] = 1.000007626290974
$arr[12,11] = 1.000007626290974
$arr[12,12] = 0.9999885685731138
$arr[12,13] = 0.9999885685731138
$arr[12,14] = 0.9999867891949984
$arr[12,15] = 0.9999949211457336
$arr[12,16] = 0.9999949211457336
$arr[12,17] = 0.9999980974320437
$arr[12,18] = 0.9999980974320437
$arr[12,19] = 0.999999323852335
$arr[13,0] = 12
$arr[13,1] = 'OffsetFTD'
$arr[13,2] = 1.002089797253162
$arr[13,3] = 0.9954024506878769
$arr[13,4] = 1.005572790015801
$arr[13,5] = 1.002089797253162
$arr[13,6] = 0.9916408172691078
$arr[13,7] = 1.002089797253162
$arr[13,8] = 0.9916408172691078
$arr[13,9] = 1.002089797253162
$arr[13,10] = 1.002089797253162
$arr[13,11] = 1.002089797253162
$arr[13,12] = 0.9968653072611348
$arr[13,13] = 0.9968653072611348
$arr[13,14] = 0.9963776884033821
$arr[13,15] = 0.9986068039251438
$arr[13,16] = 0.9986068039251438
$arr[13,17] = 0.9994775522571482
$arr[13,18] = 0.9994775522571482
$arr[13,19] = 0.9998142416220452
$arr[14,0] = 13
$arr[14,1] = 'OffsetATD'
$arr[14,2] = 1.000561146984666
$arr[14,3] = 0.9987654929139987
$arr[14,4] = 1.001496364492529
$arr[14,5] = 1.000561146984666
$arr[14,6] = 0.9977554399469895
$arr[14,7] = 1.000561146984666
$arr[14,8] = 0.9977554399469895
$arr[14,9] = 1.000561146984666
$arr[14,10] = 1.000561146984666
$arr[14,11] = 1.000561146984666
$arr[14,12] = 0.9991582934658279
$arr[14,13] = 0.9991582934658279
$arr[14,14] = 0.9990273599485514
$arr[14,15] = 0.9996259113054405
$arr[14,16] = 0.9996259113054405
$arr[14,17] = 0.9998597202252471
$arr[14,18] = 0.9998597202252471
$arr[14,19] = 0.9999501230512525
$arr[15,0] = 14
$arr[15,1] = 'Holden2.5'
$arr[15,2] = 1.000318708767534
$arr[15,3] = 0.9992988290654526
$arr[15,4] = 1.000849909411802
$arr[15,5] = 1.000318708767534
$arr[15,6] = 0.9987251366199655
$arr[15,7] = 1.000318708767534
$arr[15,8] = 0.9987251366199655
$arr[15,9] = 1.000318708767534
$arr[15,10] = 1.000318708767534
$arr[15,11] = 1.000318708767534
$arr[15,12] = 0.9995219226937497
$arr[15,13] = 0.9995219226937497
$arr[15,14] = 0.9994475581509841
$arr[15,15] = 0.9997875180516779
$arr[15,16] = 0.9997875180516779
$arr[15,17] = 0.9999203157306418
$arr[15,18] = 0.9999203157306418
$arr[15,19] = 0.9999716668999703
$arr[16,0] = 15
$arr[16,1] = 'Holden5'
$arr[16,2] = 1.000260972917643
$arr[16,3] = 0.9994258526412551
$arr[16,4] = 1.000695939983929
$arr[16,5] = 1.000260972917643
$arr[16,6] = 0.9989560899721304
$arr[16,7] = 1.000260972917643
$arr[16,8] = 0.9989560899721304
$arr[16,9] = 1.000260972917643
$arr[16,10] = 1.000260972917643
$arr[16,11] = 1.000260972917643
$arr[16,12] = 0.9996085314448868
$arr[16,13] = 0.9996085314448868
$arr[16,14] = 0.9995476385103429
$arr[16,15] = 0.9998260119358057
$arr[16,16] = 0.9998260119358057
$arr[16,17] = 0.999934752181265
$arr[16,18] = 0.999934752181265
$arr[16,19] = 0.9999768002250407
$arr[17,0] = 16
$arr[17,1] = 'Holden10'
$arr[17,2] = 1.000145233265121
$arr[17,3] = 0.9996804886355188
$arr[17,4] = 1.000387289748782
$arr[17,5] = 1.000145233265121
$arr[17,6] = 0.9994190665781356
$arr[17,7] = 1.000145233265121
$arr[17,8] = 0.9994190665781356
$arr[17,9] = 1.000145233265121
$arr[17,10] = 1.000145233265121
$arr[17,11] = 1.000145233265121
$arr[17,12] = 0.9997821499216282
$arr[17,13] = 0.9997821499216282
$arr[17,14] = 0.9997482628262584
$arr[17,15] = 0.9999031777027924
$arr[17,16] = 0.9999031777027924
$arr[17,17] = 0.9999636915933745
$arr[17,18] = 0.9999636915933745
$arr[17,19] = 0.9999870907929664
$arr[18,0] = 17
$arr[18,1] = 'Holden15'
$arr[18,2] = 1.000149368771602
$arr[18,3] = 0.9996713865245302
$arr[18,4] = 1.000398322385714
$arr[18,5] = 1.000149368771602
$arr[18,6] = 0.999402519203732
$arr[18,7] = 1.000149368771602
$arr[18,8] = 0.999402519203732
$arr[18,9] = 1.000149368771602
$arr[18,10] = 1.000149368771602
$arr[18,11] = 1.000149368771602
$arr[18,12] = 0.999775943987667
$arr[18,13] = 0.999775943987667
$arr[18,14] = 0.9997410914999548
$arr[18,15] = 0.9999004189156452
$arr[18,16] = 0.9999004189156452
$arr[18,17] = 0.9999626563796344
$arr[18,18] = 0.9999626563796344
$arr[18,19] = 0.999986722404797
$arr[19,0] = 18
$arr[19,1] = 'HexGrid-90degTilt2.5degRes'
$arr[19,2] = 1.000000068748989
$arr[19,3] = 0.9999998524974475
$arr[19,4] = 1.000000179676352
$arr[19,5] = 1.000000068748989
$arr[19,6] = 0.9999997290767715
$arr[19,7] = 1.000000068748989
$arr[19,8] = 0.9999997290767715
$arr[19,9] = 1.000000068748989
$arr[19,10] = 1.000000068748989
$arr[19,11] = 1.000000068748989
$arr[19,12] = 0.9999998989128803
$arr[19,13] = 0.9999998989128803
$arr[19,14] = 0.9999998834410694
$arr[19,15] = 0.9999999555249165
$arr[19,16] = 0.9999999555249165
$arr[19,17] = 0.9999999838309347
$arr[19,18] = 0.9999999838309347
$arr[19,19] = 0.9999999945829231
$arr[20,0] = 19
$arr[20,1] = 'HexGrid-90degTilt5degRes'
$arr[20,2] = 1.00000142974139
$arr[20,3] = 0.9999968575015092
$arr[20,4] = 1.000003810067097
$arr[20,5] = 1.00000142974139
$arr[20,6] = 0.999994283592983
$arr[20,7] = 1.00000142974139
$arr[20,8] = 0.999994283592983
$arr[20,9] = 1.00000142974139
$arr[20,10] = 1.00000142974139
$arr[20,11] = 1.00000142974139
$arr[20,12] = 0.9999978566671865
$arr[20,13] = 0.9999978566671865
$arr[20,14] = 0.9999975236119608
$arr[20,15] = 0.9999990476919211
$arr[20,16] = 0.9999990476919211
$arr[20,17] = 0.9999996432042884
$arr[20,18] = 0.9999996432042884
$arr[20,19] = 0.9999998733976266
$arr[21,0] = 20
$arr[21,1] = 'HexGrid-90degTilt10degRes'
$arr[21,2] = 1.000005027433818
$arr[21,3] = 0.9999889348666909
$arr[21,4] = 1.000013414743481
$arr[21,5] = 1.000005027433818
$arr[21,6] = 0.9999798803156478
$arr[21,7] = 1.000005027433818
$arr[21,8] = 0.9999798803156478
$arr[21,9] = 1.000005027433818
$arr[21,10] = 1.000005027433818
$arr[21,11] = 1.000005027433818
$arr[21,12] = 0.9999924538747331
$arr[21,13] = 0.9999924538747331
$arr[21,14] = 0.9999912808720524
$arr[21,15] = 0.999996645061095
$arr[21,16] = 0.999996645061095
$arr[21,17] = 0.9999987406542759
$arr[21,18] = 0.9999987406542759
$arr[21,19] = 0.9999995520378793
$arr[22,0] = 21
$arr[22,1] = 'HexGrid-90degTilt15degRes'
$arr[22,2] = 1.000011784194715
$arr[22,3] = 0.9999740771842739
$arr[22,4] = 1.000031427864921
$arr[22,5] = 1.000011784194715
$arr[22,6] = 0.9999528592053132
$arr[22,7] = 1.000011784194715
$arr[22,8] = 0.9999528592053132
$arr[22,9] = 1.000011784194715
$arr[22,10] = 1.000011784194715
$arr[22,11] = 1.000011784194715
$arr[22,12] = 0.9999823217000141
$arr[22,13] = 0.9999823217000141
$arr[22,14] = 0.9999795735281007
$arr[22,15] = 0.999992142531581
$arr[22,16] = 0.999992142531581
$arr[22,17] = 0.9999970529473645
$arr[22,18] = 0.9999970529473645
$arr[22,19] = 0.9999989528064422
$ws.Range("A1:T23").Value = $arr

# New rows 20:23 need the same bold/bordered/centered style that column A
# and row 1 use elsewhere in the table (style index carried by A19).
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
